# Update dSF (column F) values per repull/push of data and mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -1
$ws.Range("F11").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("F22").Value = -4
$ws.Range("F23").Value = 4
$ws.Range("F25").Value = -4
